$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Consolidate BOM rows ------------------------------------------------
# The capacitor references C2, C6, C8, C5, C11 get folded into their sibling
# rows (C1, C3, C9 respectively), so those rows are removed. Delete from the
# bottom up so earlier row numbers stay valid.
# Row map (before deletion):
#   3 C1   4 C10  5 C11  6 C2   7 C3   8 C4   9 C5  10 C6  11 C7  12 C8  13 C9
#  14 D1  15 J6  16 R1  17 R4  18 RV1 19 U1  20 U2  21 X1
$ws.Rows.Item(12).Delete()   # C8  -> merged into C1 row
$ws.Rows.Item(10).Delete()   # C6  -> merged into C1 row
$ws.Rows.Item(9).Delete()    # C5  -> merged into C3 row
$ws.Rows.Item(6).Delete()    # C2  -> merged into C1 row
$ws.Rows.Item(5).Delete()    # C11 -> merged into C9 row

# Row map (after deletion):
#   3 C1   4 C10  5 C3   6 C4   7 C7   8 C9   9 D1  10 J6  11 R1  12 R4
#  13 RV1 14 U1  15 U2  16 X1

# --- 2. Update consolidated reference designators / text -------------------
$ws.Range("A3").Value = "C1, C2, C6, C8"
$ws.Range("A5").Value = "C3, C5"
$ws.Range("A8").Value = "C9, C11"
$ws.Range("C8").Value = "capacitor, tantalum"

# --- 3. Update the LED row (D1): description, footprint, DigiKey link ------
$ws.Range("C9").Value = "LED, reverse mount"
$ws.Range("D9").Value = "'1206"
$ws.Range("E9").Value = "https://www.digikey.com/products/en?keywords=1497-1377-1-ND"

# --- 4. Update R4 value ------------------------------------------------------
$ws.Range("B12").Value = "25.5K"

# --- 5. Turn the DigiKey column links into real hyperlinks ------------------
$ws.Hyperlinks.Add($ws.Range("E9"), $ws.Range("E9").Text) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E13"), $ws.Range("E13").Text) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E14"), $ws.Range("E14").Text) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E15"), $ws.Range("E15").Text) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E16"), $ws.Range("E16").Text) | Out-Null

# --- 6. Selection / active cell ---------------------------------------------
$ws.Range("E16").Select()
